$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column A (a bare numeric index 0..19, styled like the header)
# is no longer wanted: the dataset-name labels that used to live in
# column B should become the new column A, and everything to the right
# shifts left by one. Deleting the entire column A and letting Excel
# shift the remaining columns left accomplishes exactly that.
$ws.Columns.Item(1).Delete()
